$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 40.942832
$ws.Range("H2").Value = 122.828496
$ws.Range("I2").Value = 0.2583000005785167
$ws.Range("J2").Value = 0.2583000005785167
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 87.94127933333334
$ws.Range("N2").Value = 263.823838
$ws.Range("O2").Value = 0.4109331243514438
$ws.Range("P2").Value = 0.4109331243514437
$ws.Range("Q2").Value = 3600.565025609739
$ws.Range("R2").Value = 32405.08523048765
$ws.Range("S2").Value = 0.1061440262577096
$ws.Range("T2").Value = 0.1061440262577096

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 40.942832
$ws.Range("H3").Value = 122.828496
$ws.Range("I3").Value = 0.2583000005785167
$ws.Range("J3").Value = 0.2583000005785167
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 52.441971
$ws.Range("N3").Value = 157.325913
$ws.Range("O3").Value = 0.2450515065683088
$ws.Range("P3").Value = 0.2450515065683087
$ws.Range("Q3").Value = 2147.122808401872
$ws.Range("R3").Value = 19324.10527561685
$ws.Range("S3").Value = 0.06329680428836053
$ws.Range("T3").Value = 0.06329680428836053

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 40.942832
$ws.Range("H4").Value = 122.828496
$ws.Range("I4").Value = 0.2583000005785167
$ws.Range("J4").Value = 0.2583000005785167
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 54.667459
$ws.Range("N4").Value = 164.002377
$ws.Range("O4").Value = 0.255450795093328
$ws.Range("P4").Value = 0.255450795093328
$ws.Range("Q4").Value = 2238.240589703888
$ws.Range("R4").Value = 20144.16530733499
$ws.Range("S4").Value = 0.06598294052038917
$ws.Range("T4").Value = 0.06598294052038915

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 40.942832
$ws.Range("H5").Value = 122.828496
$ws.Range("I5").Value = 0.2583000005785167
$ws.Range("J5").Value = 0.2583000005785167
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 18.95316166666667
$ws.Range("N5").Value = 56.85948500000001
$ws.Range("O5").Value = 0.08856457398691947
$ws.Range("P5").Value = 0.08856457398691944
$ws.Range("Q5").Value = 775.9961139871735
$ws.Range("R5").Value = 6983.96502588456
$ws.Range("S5").Value = 0.02287622951205738
$ws.Range("T5").Value = 0.02287622951205737

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 79.68771233333334
$ws.Range("H6").Value = 239.063137
$ws.Range("I6").Value = 0.5027335710876245
$ws.Range("J6").Value = 0.5027335710876245
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 87.94127933333334
$ws.Range("N6").Value = 263.823838
$ws.Range("O6").Value = 0.4109331243514438
$ws.Range("P6").Value = 0.4109331243514437
$ws.Range("Q6").Value = 7007.83936973998
$ws.Range("R6").Value = 63070.55432765982
$ws.Range("S6").Value = 0.2065898770833962
$ws.Range("T6").Value = 0.2065898770833962

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 79.68771233333334
$ws.Range("H7").Value = 239.063137
$ws.Range("I7").Value = 0.5027335710876245
$ws.Range("J7").Value = 0.5027335710876245
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 52.441971
$ws.Range("N7").Value = 157.325913
$ws.Range("O7").Value = 0.2450515065683088
$ws.Range("P7").Value = 0.2450515065683087
$ws.Range("Q7").Value = 4178.980699241009
$ws.Range("R7").Value = 37610.82629316908
$ws.Range("S7").Value = 0.1231956189974883
$ws.Range("T7").Value = 0.1231956189974883

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 79.68771233333334
$ws.Range("H8").Value = 239.063137
$ws.Range("I8").Value = 0.5027335710876245
$ws.Range("J8").Value = 0.5027335710876245
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 54.667459
$ws.Range("N8").Value = 164.002377
$ws.Range("O8").Value = 0.255450795093328
$ws.Range("P8").Value = 0.255450795093328
$ws.Range("Q8").Value = 4356.324746786295
$ws.Range("R8").Value = 39206.92272107665
$ws.Range("S8").Value = 0.1284236904544418
$ws.Range("T8").Value = 0.1284236904544418

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 79.68771233333334
$ws.Range("H9").Value = 239.063137
$ws.Range("I9").Value = 0.5027335710876245
$ws.Range("J9").Value = 0.5027335710876245
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 18.95316166666667
$ws.Range("N9").Value = 56.85948500000001
$ws.Range("O9").Value = 0.08856457398691947
$ws.Range("P9").Value = 0.08856457398691944
$ws.Range("Q9").Value = 1510.334094700494
$ws.Range("R9").Value = 13593.00685230445
$ws.Range("S9").Value = 0.04452438455229816
$ws.Range("T9").Value = 0.04452438455229814

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 37.53186833333334
$ws.Range("H10").Value = 112.595605
$ws.Range("I10").Value = 0.2367809244903433
$ws.Range("J10").Value = 0.2367809244903433
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 87.94127933333334
$ws.Range("N10").Value = 263.823838
$ws.Range("O10").Value = 0.4109331243514438
$ws.Range("P10").Value = 0.4109331243514437
$ws.Range("Q10").Value = 3300.600517003555
$ws.Range("R10").Value = 29705.40465303199
$ws.Range("S10").Value = 0.09730112508764006
$ws.Range("T10").Value = 0.09730112508764004

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 37.53186833333334
$ws.Range("H11").Value = 112.595605
$ws.Range("I11").Value = 0.2367809244903433
$ws.Range("J11").Value = 0.2367809244903433
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 52.441971
$ws.Range("N11").Value = 157.325913
$ws.Range("O11").Value = 0.2450515065683088
$ws.Range("P11").Value = 0.2450515065683087
$ws.Range("Q11").Value = 1968.245150712485
$ws.Range("R11").Value = 17714.20635641237
$ws.Range("S11").Value = 0.05802352227299558
$ws.Range("T11").Value = 0.05802352227299557

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 37.53186833333334
$ws.Range("H12").Value = 112.595605
$ws.Range("I12").Value = 0.2367809244903433
$ws.Range("J12").Value = 0.2367809244903433
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 54.667459
$ws.Range("N12").Value = 164.002377
$ws.Range("O12").Value = 0.255450795093328
$ws.Range("P12").Value = 0.255450795093328
$ws.Range("Q12").Value = 2051.771873305899
$ws.Range("R12").Value = 18465.94685975309
$ws.Range("S12").Value = 0.06048587542399146
$ws.Range("T12").Value = 0.06048587542399144

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 37.53186833333334
$ws.Range("H13").Value = 112.595605
$ws.Range("I13").Value = 0.2367809244903433
$ws.Range("J13").Value = 0.2367809244903433
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 18.95316166666667
$ws.Range("N13").Value = 56.85948500000001
$ws.Range("O13").Value = 0.08856457398691947
$ws.Range("P13").Value = 0.08856457398691944
$ws.Range("Q13").Value = 711.3475681737141
$ws.Range("R13").Value = 6402.128113563426
$ws.Range("S13").Value = 0.0209704017057162
$ws.Range("T13").Value = 0.02097040170571619

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.3464216666666666
$ws.Range("H14").Value = 1.039265
$ws.Range("I14").Value = 0.002185503843515531
$ws.Range("J14").Value = 0.002185503843515531
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 87.94127933333334
$ws.Range("N14").Value = 263.823838
$ws.Range("O14").Value = 0.4109331243514438
$ws.Range("P14").Value = 0.4109331243514437
$ws.Range("Q14").Value = 30.46476455545222
$ws.Range("R14").Value = 274.18288099907
$ws.Range("S14").Value = 0.0008980959226979262
$ws.Range("T14").Value = 0.0008980959226979261

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.3464216666666666
$ws.Range("H15").Value = 1.039265
$ws.Range("I15").Value = 0.002185503843515531
$ws.Range("J15").Value = 0.002185503843515531
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 52.441971
$ws.Range("N15").Value = 157.325913
$ws.Range("O15").Value = 0.2450515065683088
$ws.Range("P15").Value = 0.2450515065683087
$ws.Range("Q15").Value = 18.167034997105
$ws.Range("R15").Value = 163.503314973945
$ws.Range("S15").Value = 0.0005355610094643102
$ws.Range("T15").Value = 0.0005355610094643102

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.3464216666666666
$ws.Range("H16").Value = 1.039265
$ws.Range("I16").Value = 0.002185503843515531
$ws.Range("J16").Value = 0.002185503843515531
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 54.667459
$ws.Range("N16").Value = 164.002377
$ws.Range("O16").Value = 0.255450795093328
$ws.Range("P16").Value = 0.255450795093328
$ws.Range("Q16").Value = 18.93799225921166
$ws.Range("R16").Value = 170.441930332905
$ws.Range("S16").Value = 0.0005582886945055669
$ws.Range("T16").Value = 0.0005582886945055668

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.3464216666666666
$ws.Range("H17").Value = 1.039265
$ws.Range("I17").Value = 0.002185503843515531
$ws.Range("J17").Value = 0.002185503843515531
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 18.95316166666667
$ws.Range("N17").Value = 56.85948500000001
$ws.Range("O17").Value = 0.08856457398691947
$ws.Range("P17").Value = 0.08856457398691944
$ws.Range("Q17").Value = 6.565785853169445
$ws.Range("R17").Value = 59.092072678525
$ws.Range("S17").Value = 0.0001935582168477281
$ws.Range("T17").Value = 0.0001935582168477281
